$d = $word.ActiveDocument

# --- First paragraph ("**ID__AFFARS_5318_topic_7__ID** ") -----------------
$para = $d.Paragraphs(1)

# Replace the placeholder id text. The paragraph currently holds two runs:
# "**ID__AFFARS_5318_topic_7__ID**" followed by a run containing a single
# trailing space. Searching/replacing across the whole paragraph range
# (including that trailing space) merges them into a single run holding
# just the new id, dropping the now-unwanted space-only run.
$para.Range.Find.Execute("**ID__AFFARS_5318_topic_7__ID** ", $false, $false, $false, $false, $false, `
                          $true, 1, $false, "**ID__AFFARS_5318_202__ID**", 2)

# Give the paragraph the same "boxed" look already used by the "See ..."
# paragraph further down: a border with 5-twip spacing on every side and a
# deeper left indent (120 -> 225 twips).
$para.Format.Borders.DistanceFromTop = 5
$para.Format.Borders.DistanceFromLeft = 5
$para.Format.Borders.DistanceFromBottom = 5
$para.Format.Borders.DistanceFromRight = 5

# LeftIndent is expressed in points by the object model; 225 twips is
# 11.25 points (20 twips per point).
$para.Format.LeftIndent = 11.25
